# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update A30: "John" -> "John-analyst"
$ws.Range("A30").Value = "John-analyst"

# Update D39/E39 with new title/link
$ws.Range("D39").Value = "A Gentle Introduction to Chefboost for Applied Machine Learning"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/A-Gentle-Introduction-to-Chefboost-for-Applied-Machine-Learning-1"

# Ensure F40 is set (used flag)
$ws.Range("F40").Value = "o"

# Append new row 41
$ws.Range("A41").Value = "cloudinsight"
$ws.Range("B41").Value = "http://cloudinsight.net/feed/"
$ws.Range("C41").Value = "etc"
$ws.Range("D41").Value = "머신 러닝 모델 관리법"
$ws.Range("E41").Value = "http://cloudinsight.net/ai/%eb%a8%b8%ec%8b%a0-%eb%9f%ac%eb%8b%9d%ec%9d%98-%eb%aa%a8%eb%8d%b8-%ea%b4%80%eb%a6%ac%eb%b2%95/"
$ws.Range("F41").Value = "o"
